# fix(allowance): issue on import allowance generate payroll
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix shared string typo: "Employee id" -> "Employee Id"
$ws.Range("A1").Value = "Employee Id"

# 2. Apply a Text number format (numFmtId 49, "@") to column A so employee
#    ids import without losing leading zeros / being coerced to numbers.
$ws.Columns("A").NumberFormat = "@"

# 3. Move the active selection to I12 (matches the recorded selection state).
$ws.Range("I12").Select()

# 4. Set up the page for printing (adds pageSetup entry on save).
$ws.PageSetup.Orientation = 1
